$d = $word.ActiveDocument

# Move to the very end of the document body (after "Rule 13 [IMP]")
$insertionPoint = $d.Content
$insertionPoint.Collapse(0)

# Append the new "19 December 2022" class-notes block:
#   - one blank paragraph
#   - a bold/underlined date heading paragraph
#   - a paragraph with "Sec 50: Rule of Damdupat" (Damdupat flagged as a
#     spell-check exception via proofErr, matching the source markup)
$newParagraphsXml = @'
<w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="24"/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>19 December 2022</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Sec 50: Rule of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Damdupat</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@

[void]$insertionPoint.InsertXML($newParagraphsXml)
